# 🚌 141: 31/12 18:00 LP1912+6203+6173
#
# Appends freshly-scraped rows to the "LP1912" and "6203-6173" sheets and
# refreshes the "Última actualización" / "Total filas" banner cells on all
# three sheets (the scraper re-stamps every sheet's timestamp even when a
# given sheet gets no new rows, as on "LP1912-215" here).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$stamp = "Última actualización: 31/12/2025 15:00:03"

# --- Sheet "LP1912": 15 new rows (1026-1040) ------------------------------
# Columns: B=Hora_Scrap, C=Hora_Llegada, D=Línea, E=Minutos(n), F=Parada, G=Fecha
$sheet1Rows = @(
    @(1026, '14:59:52', '15:01', '14_ABASTO', 2, 'LP1912', '31/12/2025'),
    @(1027, '14:59:52', '15:03', '23_HERNANDEZ', 4, 'LP1912', '31/12/2025'),
    @(1028, '14:59:52', '15:13', '15_ABASTO', 14, 'LP1912', '31/12/2025'),
    @(1029, '14:59:52', '15:24', '11_ETCHEVERRY', 25, 'LP1912', '31/12/2025'),
    @(1030, '14:59:52', '15:33', '16_SANTA ANA', 34, 'LP1912', '31/12/2025'),
    @(1031, '14:59:52', '15:33', '23_HERNANDEZ', 34, 'LP1912', '31/12/2025'),
    @(1032, '14:59:52', '15:44', '14_ABASTO', 45, 'LP1912', '31/12/2025'),
    @(1033, '14:59:52', '15:57', '16_SANTA ANA', 58, 'LP1912', '31/12/2025'),
    @(1034, '14:59:52', '16:01', '15_ABASTO', 62, 'LP1912', '31/12/2025'),
    @(1035, '14:59:52', '16:09', '16_SANTA ANA', 70, 'LP1912', '31/12/2025'),
    @(1036, '14:59:52', '16:16', '10_OLMOS', 77, 'LP1912', '31/12/2025'),
    @(1037, '14:59:52', '16:24', '11_ETCHEVERRY', 85, 'LP1912', '31/12/2025'),
    @(1038, '14:59:52', '16:31', '16_P MOR-SANTA ANA', 92, 'LP1912', '31/12/2025'),
    @(1039, '14:59:52', '16:33', '23_HERNANDEZ', 94, 'LP1912', '31/12/2025'),
    @(1040, '14:59:52', '16:34', '225_GOMEZ', 95, 'LP1912', '31/12/2025')
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

$ws1.Range("A2").Value = $stamp
$ws1.Range("A3").Value = "Total filas: 1039"

# --- Sheet "LP1912-215": timestamp refresh only, no new rows -------------
$ws2.Range("A2").Value = $stamp

# --- Sheet "6203-6173": 1 new row (129) -----------------------------------
# Columns: B=Fecha, C=Hora_Scrap, D=Hora_Llegada, E=Línea, F=Minutos(n), G=Parada
$ws3.Cells.Item(129, 2).Value = "31/12/2025"
$ws3.Cells.Item(129, 3).Value = "14:59:58"
$ws3.Cells.Item(129, 4).Value = "15:46"
$ws3.Cells.Item(129, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(129, 6).Value = 47
$ws3.Cells.Item(129, 7).Value = "L6203"

$ws3.Range("A2").Value = $stamp
$ws3.Range("A3").Value = "Total filas: 128"
